# Generate Report for Archive
# - Update the localization "Status" text from "Ready for handoff" to
#   "In Translation" on every sheet that references it, and shrink the
#   now-narrower Status column(s) to match (mirrors an Excel AutoFit/resize
#   after the shorter text was entered).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Closest value reachable through ColumnWidth's character-unit rounding to
# the target OOXML column width (~13.4101845877511).
$newColWidth = 12.5

# --- Overview sheet: Status appears twice, in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet: Status is column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet: Status is column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
